# Apply the "02_clean transform" changes:
# - Add two new lookup rows (grao, coco) distinguishing unit_from/unit_to
#   from the unidade_padrao ("_toneladas") reference values used previously.
# - Update D34/E34 (cafe, 1974-2001) from "coco_toneladas" to "coco"
# - Update D35/E35 (cafe, 2002-9999) from "grao_toneladas" to "grao"
# - Update the active sheet view (scroll position / selection) to reflect
#   where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update unit_from / unit_to values for the cafe rows (row 34 and row 35).
# Order matters for shared-string table layout: "grao" is appended before
# "coco" so the new entries land at indices 35 ("grao") and 36 ("coco").
$ws.Range("D35").Value = "grao"
$ws.Range("E35").Value = "grao"

$ws.Range("D34").Value = "coco"
$ws.Range("E34").Value = "coco"

# Reflect the final view/selection state (scrolled down, active cell E34)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E34").Select()
